# Generate Report for Handoff
# Adds a new row (row 3) to each of the three sheets (Overview, zh-cn, de-de)
# reflecting a newly-generated handoff package for file
# d62cc0f4-442f-4868-b32b-a53a6dd3caf5...md ("Ready for handoff").

$wb = $excel.ActiveWorkbook

$newMdFile    = 'd62cc0f4-442f-4868-b32b-a53a6dd3caf5ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md'
$newMdDisplay = 'e2e\d62cc0f4-442f-4868-b32b-a53a6dd3caf5ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md'
$statusReady  = 'Ready for handoff'
$hoDateOview  = '2016-08-13 20:35:39'
$xlfZhCn      = 'd62cc0f4-442f-4868-b32b-a53a6dd3caf5ooooooooooooooooooooooooooooooooooooooooooo.9dbfaa09639159c708659481fa686792228797cd.zh-cn.xlf'
$xlfZhCnDate  = '2016-08-13 20:35:31'
$xlfDeDe      = 'd62cc0f4-442f-4868-b32b-a53a6dd3caf5ooooooooooooooooooooooooooooooooooooooooooo.9dbfaa09639159c708659481fa686792228797cd.de-de.xlf'
$hyperlinkUrl = 'https://github.com/OpenLocalizationTestOrg/oltest/blob/8b2a21de809576e2d18afc45180465c043f211ef/e2e/d62cc0f4-442f-4868-b32b-a53a6dd3caf5ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md'

$dtFormat = 'yyyy-mm-dd HH:mm:ss'
$hlColor  = 15570276   # RGB(0x64,0x95,0xED) -> matches existing HyperLink font color FF6495ED

function Set-DateTimeStyle($range) {
    $range.NumberFormat = $dtFormat
}

function Set-HyperlinkFontStyle($range) {
    $range.Font.Underline = $true
    $range.Font.Color = $hlColor
}

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A3").Value = $newMdFile
$wsOverview.Range("B3").Value = $newMdDisplay
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("D3").Value = ""
$wsOverview.Range("E3").Value = $statusReady
$wsOverview.Range("F3").Value = $statusReady
$wsOverview.Range("G3").Value = $hoDateOview

Set-HyperlinkFontStyle $wsOverview.Range("B3")
Set-DateTimeStyle $wsOverview.Range("G3")

$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), $hyperlinkUrl, "", "", $newMdDisplay) | Out-Null

$wsOverview.Columns.Item(5).AutoFit() | Out-Null
$wsOverview.Columns.Item(6).AutoFit() | Out-Null

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.ListRows.Add() | Out-Null

$wsZhCn.Range("A3").Value = $newMdFile
$wsZhCn.Range("B3").Value = ".md"
$wsZhCn.Range("C3").Value = $statusReady
$wsZhCn.Range("D3").Value = "e2e"
$wsZhCn.Range("E3").Value = "ht"
$wsZhCn.Range("F3").Value = "False"
$wsZhCn.Range("G3").Value = $xlfZhCn
$wsZhCn.Range("H3").Value = $xlfZhCnDate
$wsZhCn.Range("I3").Value = ""
$wsZhCn.Range("J3").Value = ""
$wsZhCn.Range("K3").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("L3").Value = ""
$wsZhCn.Range("M3").Value = "True"
$wsZhCn.Range("N3").Value = ""
$wsZhCn.Range("O3").Value = "False"
$wsZhCn.Range("P3").Value = ""

Set-HyperlinkFontStyle $wsZhCn.Range("A3")
Set-DateTimeStyle $wsZhCn.Range("H3")
Set-DateTimeStyle $wsZhCn.Range("K3")

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $hyperlinkUrl, "", "", $newMdFile) | Out-Null

$wsZhCn.Columns.Item(3).AutoFit() | Out-Null

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.ListRows.Add() | Out-Null

$wsDeDe.Range("A3").Value = $newMdFile
$wsDeDe.Range("B3").Value = ".md"
$wsDeDe.Range("C3").Value = $statusReady
$wsDeDe.Range("D3").Value = "e2e"
$wsDeDe.Range("E3").Value = "ht"
$wsDeDe.Range("F3").Value = "False"
$wsDeDe.Range("G3").Value = $xlfDeDe
$wsDeDe.Range("H3").Value = $hoDateOview
$wsDeDe.Range("I3").Value = ""
$wsDeDe.Range("J3").Value = ""
$wsDeDe.Range("K3").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("L3").Value = ""
$wsDeDe.Range("M3").Value = "True"
$wsDeDe.Range("N3").Value = ""
$wsDeDe.Range("O3").Value = "False"
$wsDeDe.Range("P3").Value = ""

Set-HyperlinkFontStyle $wsDeDe.Range("A3")
Set-DateTimeStyle $wsDeDe.Range("H3")
Set-DateTimeStyle $wsDeDe.Range("K3")

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $hyperlinkUrl, "", "", $newMdFile) | Out-Null

$wsDeDe.Columns.Item(3).AutoFit() | Out-Null
